# "add kakao logup - login"
# Adds a new ErrorMessages row for the kakao login module and updates the
# active selections left behind by the edit session (Status!B4, then
# ErrorMessages!D32 on the newly-added row).

$wb = $excel.ActiveWorkbook

# --- Status sheet: user last left the selection on B4 -----------------
$wsStatus = $wb.Worksheets.Item("Status")
$wsStatus.Activate()
[void]$wsStatus.Range("B4").Select()

# --- ErrorMessages sheet: append a new row for the kakao logIn module -
$wsErrors = $wb.Worksheets.Item("ErrorMessages")
$wsErrors.Activate()

$wsErrors.Cells.Item(32, 1).Value = 30
$wsErrors.Cells.Item(32, 2).Value = "modulse/user.User.logIn"
$wsErrors.Cells.Item(32, 3).Value = 400
$wsErrors.Cells.Item(32, 4).Value = "deleted"

# Match the formatting used by the row above (B/D columns carry style 1)
$wsErrors.Range("B31:D31").Copy()
$wsErrors.Range("B32:D32").PasteSpecial(-4122)

# Leave the selection where the editor ended up: the new D32 cell
[void]$wsErrors.Range("D32").Select()
